$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Cells.Item(2, 4)
$cell.NumberFormat = "@"
$cell.Value = "43.703.02"
$cell.Style = "Normal"
$ws.Cells.Item(2, 5).Value = "  +1.81%  "
$cell = $ws.Cells.Item(3, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.339.73"
$cell.Style = "Normal"
$ws.Cells.Item(3, 5).Value = "  +2.14%  "
$ws.Cells.Item(4, 5).Value = "  +0.02%  "
$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = "312.16"
$cell.Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  -0.69%  "
$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = "108.13"
$cell.Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  +2.77%  "
$ws.Cells.Item(7, 5).Value = "  +0.88%  "
$ws.Cells.Item(8, 5).Value = "  +0.04%  "
$cell = $ws.Cells.Item(9, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.619"
$cell.Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  +2.18%  "
$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = "@"
$cell.Value = "41.23"
$cell.Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  +3.74%  "
$cell = $ws.Cells.Item(11, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0918"
$cell.Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  +1.26%  "
$cell = $ws.Cells.Item(12, 4)
$cell.NumberFormat = "@"
$cell.Value = "8.55"
$cell.Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  +1.48%  "
$ws.Cells.Item(13, 5).Value = "  -0.91%  "
$ws.Cells.Item(14, 5).Value = "  +0.94%  "
$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = "@"
$cell.Value = "15.49"
$cell.Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  +1.37%  "
$cell = $ws.Cells.Item(16, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.696.63"
$cell.Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  +2.19%  "
$cell = $ws.Cells.Item(17, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.331.51"
$cell.Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  -0.60%  "
$cell = $ws.Cells.Item(18, 4)
$cell.NumberFormat = "@"
$cell.Value = "43.858.76"
$cell.Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  +2.39%  "
$ws.Cells.Item(19, 5).Value = "  +1.85%  "
$ws.Cells.Item(20, 5).Value = "  +1.21%  "
$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = "@"
$cell.Value = "12.97"
$cell.Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  -5.74%  "
$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = "@"
$cell.Value = "74.17"
$cell.Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  +0.26%  "
$ws.Cells.Item(23, 5).Value = "  -2.78%  "
$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = "@"
$cell.Value = "268.67"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.27"
$cell.Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  +2.65%  "
$ws.Cells.Item(26, 5).Value = "  -0.15%  "
$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = "@"
$cell.Value = "7.61"
$cell.Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  +2.67%  "
$cell = $ws.Cells.Item(28, 4)
$cell.NumberFormat = "@"
$cell.Value = "11.11"
$cell.Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  +2.06%  "
$ws.Cells.Item(29, 5).Value = "  -2.00%  "
$cell = $ws.Cells.Item(30, 4)
$cell.NumberFormat = "@"
$cell.Value = "39.02"
$cell.Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  +4.40%  "
$cell = $ws.Cells.Item(31, 4)
$cell.NumberFormat = "@"
$cell.Value = "22.62"
$cell.Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  +0.41%  "
$cell = $ws.Cells.Item(32, 4)
$cell.NumberFormat = "@"
$cell.Value = "168.88"
$cell.Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  +0.99%  "
$ws.Cells.Item(33, 5).Value = "  +1.15%  "
$cell = $ws.Cells.Item(34, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.84"
$cell.Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  +9.95%  "
$ws.Cells.Item(35, 5).Value = "  +0.64%  "
$cell = $ws.Cells.Item(36, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.115"
$cell.Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  +1.08%  "
$cell = $ws.Cells.Item(37, 4)
$cell.NumberFormat = "@"
$cell.Value = "4.73"
$cell.Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  +3.75%  "
$ws.Cells.Item(38, 5).Value = "  +3.14%  "
$ws.Cells.Item(39, 5).Value = "  +9.61%  "
$ws.Cells.Item(40, 5).Value = "  -0.93%  "
$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.72"
$cell.Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  +8.83%  "
$ws.Cells.Item(42, 5).Value = "  +10.56%  "
$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.238"
$cell.Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  +2.33%  "
$ws.Cells.Item(44, 5).Value = "  +9.88%  "
$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = "@"
$cell.Value = "71.51"
$cell.Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  +0.87%  "
$ws.Cells.Item(46, 5).Value = "  -0.02%  "
$cell = $ws.Cells.Item(47, 4)
$cell.NumberFormat = "@"
$cell.Value = "113.70"
$cell.Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  -0.11%  "
$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.675.31"
$cell.Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  -3.69%  "
$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = "@"
$cell.Value = "76.98"
$cell.Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  -3.70%  "
$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = "@"
$cell.Value = "8.93"
$cell.Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  +1.93%  "
$ws.Cells.Item(51, 5).Value = "  +13.56%  "
